{"js": "// Update the \"Project Scarlett\" codename reference to the public Xbox\n// Series X|S devkit naming (Xbox GDK \"November GDK release\" update).\n//\n// Before: \"If using Project Scarlett, set the active solution platform to \"\n// After:  \"If using an Xbox Series X|S devkit, set the active solution platform to \"\n//\n// The sentence is rebuilt as three runs (matching how the text reads once\n// split around the newly-inserted phrase): \"If using \", the new phrase, and\n// \", set the active solution platform to \".\n\nconst body = context.document.body;\n\n// Locate the old codename text inside the paragraph (\"If using Project\n// Scarlett, set the active solution platform to Gaming.Xbox.Scarlett.x64.\").\nconst searchResults = body.search(\"Project Scarlett\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const oldPhrase = searchResults.items[0];\n\n  // Replace \"Project Scarlett\" with the new phrasing in place.\n  oldPhrase.insertText(\"an Xbox Series X|S devkit\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Re-locate the freshly inserted text so it ends up in its own run,\n  // distinct from the surrounding \"If using \" / \", set the active...\" text\n  // (mirrors the three-run split produced when this text was originally\n  // authored).\n  const newResults = body.search(\"an Xbox Series X|S devkit\", { matchCase: true });\n  newResults.load(\"items\");\n  await context.sync();\n\n  const newPhrase = newResults.items[0];\n  newPhrase.font.bold = true;\n  await context.sync();\n  newPhrase.font.bold = false;\n  await context.sync();\n}\n", "ps1": "# Update the \"Project Scarlett\" codename reference to the public Xbox\n# Series X|S devkit naming (Xbox GDK \"November GDK release\" update).\n#\n# Before: \"If using Project Scarlett, set the active solution platform to \"\n# After:  \"If using an Xbox Series X|S devkit, set the active solution platform to \"\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Project Scarlett\")\n\nif ($found) {\n    # Replace the codename with the new phrasing in place.\n    $rng.Text = \"an Xbox Series X|S devkit\"\n\n    # Toggle (and restore) direct formatting on the replacement text so it\n    # is written out as its own run, distinct from the surrounding\n    # \"If using \" / \", set the active solution platform to \" text (mirrors\n    # the three-run split produced when this text was originally authored).\n    $rng.Bold = 1\n    $rng.Bold = 0\n}\n"}
